$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.723.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.21%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.901.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.34%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  -0.04%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D7").Value = "'0.5226"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +5.68%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3786"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.21%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07240"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.15%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'21.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.04%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.9017"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.03%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07646"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.29%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.917.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.41%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.445"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.30%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'92.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.29%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +0.01%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.000008686"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.56%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +0.04%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'27.778.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.15%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'14.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.31%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.136"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.37%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'2.164.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.36%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.83%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.612"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.38%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'152.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.68%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.867"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.91%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'2.161"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.14%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'114.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.63%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.838"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.91%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +1.54%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.187"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.01%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.835"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +4.12%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'ImmutableX"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'0.7782"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.44%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'ARBITRUM"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'1.220"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.76%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.02094"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.26%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.576"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.82%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'3.070"
$ws.Range("D38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.59%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.5544"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.77%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.05289"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.06%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'6.727"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.57%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'117.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.09%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.502"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.56%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.1516"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.31%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4809"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.42%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'10.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.11%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.01%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.611"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.30%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'66.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.84%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.92%  "
$ws.Range("E51").Style = "Normal"
